$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 76
$ws.Range("I2").Value = 233
$ws.Range("J2").Value = 999
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 276
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = 192
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 16
$ws.Range("S2").Value = 116
$ws.Range("T2").Value = 161
$ws.Range("U2").Value = 14
$ws.Range("V2").Value = 1591
$ws.Range("X2").Value = 1543
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 34
$ws.Range("AA2").Value = 10
